$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.061321398556226
$ws.Range("D2").Value = 1.067797366434812
$ws.Range("E2").Value = 1.065018583725078
$ws.Range("F2").Value = 1.076606152553309
$ws.Range("I2").Value = 1.057324253352831
$ws.Range("J2").Value = 1.066297682308172
$ws.Range("K2").Value = 1.070504265738888
$ws.Range("L2").Value = 1.06773295471182
$ws.Range("M2").Value = 1.079289642001102
$ws.Range("N2").Value = 1.025557616551882
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.062555617903231
$ws.Range("D3").Value = 1.068818223727867
$ws.Range("E3").Value = 1.066203777349063
$ws.Range("F3").Value = 1.077759103800774
$ws.Range("I3").Value = 1.057787663313956
$ws.Range("J3").Value = 1.067184540359022
$ws.Range("K3").Value = 1.071340467621829
$ws.Range("L3").Value = 1.068732540805913
$ws.Range("M3").Value = 1.080259309246273
$ws.Range("N3").Value = 1.025878551360732
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.06335367644633
$ws.Range("D4").Value = 1.069478301842035
$ws.Range("E4").Value = 1.066970440053253
$ws.Range("F4").Value = 1.078504946511711
$ws.Range("I4").Value = 1.058085985866455
$ws.Range("J4").Value = 1.067757294141358
$ws.Range("K4").Value = 1.071880450760555
$ws.Range("L4").Value = 1.06937852672855
$ws.Range("M4").Value = 1.080885968558615
$ws.Range("N4").Value = 1.026085266381966
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.063689047899703
$ws.Range("D5").Value = 1.069755684145496
$ws.Range("E5").Value = 1.067292690386455
$ws.Range("F5").Value = 1.078818454141338
$ws.Range("I5").Value = 1.058211034274145
$ws.Range("J5").Value = 1.067997817570153
$ws.Range("K5").Value = 1.072107198949938
$ws.Range("L5").Value = 1.069649906394468
$ws.Range("M5").Value = 1.081149230985001
$ws.Range("N5").Value = 1.02617194213985
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.063745350557639
$ws.Range("D6").Value = 1.069802251182785
$ws.Range("E6").Value = 1.067346794488543
$ws.Range("F6").Value = 1.078871090884894
$ws.Range("I6").Value = 1.058232008995948
$ws.Range("J6").Value = 1.068038187183099
$ws.Range("K6").Value = 1.07214525572387
$ws.Range("L6").Value = 1.069695460981325
$ws.Range("M6").Value = 1.081193423095263
$ws.Range("N6").Value = 1.026186482070143
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.063358158212874
$ws.Range("D7").Value = 1.069482008686086
$ws.Range("E7").Value = 1.066974746190246
$ws.Range("F7").Value = 1.07850913578923
$ws.Range("I7").Value = 1.058087658207812
$ws.Range("J7").Value = 1.067760509055541
$ws.Range("K7").Value = 1.071883481603525
$ws.Range("L7").Value = 1.069382153672255
$ws.Range("M7").Value = 1.080889487009499
$ws.Range("N7").Value = 1.026086425439943
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.061738627114126
$ws.Range("D8").Value = 1.068142471243856
$ws.Range("E8").Value = 1.06541917571674
$ws.Range("F8").Value = 1.076995839022416
$ws.Range("I8").Value = 1.057481183294112
$ws.Range("J8").Value = 1.066597629493257
$ws.Range("K8").Value = 1.070787091803525
$ws.Range("L8").Value = 1.068070938986928
$ws.Range("M8").Value = 1.079617508128485
$ws.Range("N8").Value = 1.025666275236987
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.058880369064507
$ws.Range("D9").Value = 1.065778253012542
$ws.Range("E9").Value = 1.062676156854685
$ws.Range("F9").Value = 1.074327652928761
$ws.Range("I9").Value = 1.05640069766916
$ws.Range("J9").Value = 1.064539964504406
$ws.Range("K9").Value = 1.068846655178272
$ws.Range("L9").Value = 1.06575410063428
$ws.Range("M9").Value = 1.077370069277358
$ws.Range("N9").Value = 1.02491860898793
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.056971712576747
$ws.Range("D10").Value = 1.06419946576034
$ws.Range("E10").Value = 1.060846067042502
$ws.Range("F10").Value = 1.072547697459011
$ws.Range("I10").Value = 1.055672375733216
$ws.Range("J10").Value = 1.063162347510016
$ws.Range("K10").Value = 1.06754724676138
$ws.Range("L10").Value = 1.064205180694989
$ws.Range("M10").Value = 1.075867603815573
$ws.Range("N10").Value = 1.024415219083896
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.056144457882102
$ws.Range("D11").Value = 1.063515184093164
$ws.Range("E11").Value = 1.060053253765578
$ws.Range("F11").Value = 1.071776654833945
$ws.Range("I11").Value = 1.055355092342766
$ws.Range("J11").Value = 1.062564413746279
$ws.Range("K11").Value = 1.066983194346431
$ws.Range("L11").Value = 1.063533419915027
$ws.Range("M11").Value = 1.075216005745211
$ws.Range("N11").Value = 1.024196065117877
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.055837056026608
$ws.Range("D12").Value = 1.06326091074394
$ws.Range("E12").Value = 1.05975870903227
$ws.Range("F12").Value = 1.071490206336411
$ws.Range("I12").Value = 1.055236950076583
$ws.Range("J12").Value = 1.062342099450947
$ws.Range("K12").Value = 1.066773467816416
$ws.Range("L12").Value = 1.063283735146017
$ws.Range("M12").Value = 1.074973817695775
$ws.Range("N12").Value = 1.024114483174431
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.055903000371405
$ws.Range("D13").Value = 1.063315457829193
$ws.Range("E13").Value = 1.059821892570986
$ws.Range("F13").Value = 1.071551652763594
$ws.Range("I13").Value = 1.055262305097323
$ws.Range("J13").Value = 1.062389796390762
$ws.Range("K13").Value = 1.066818464510695
$ws.Range("L13").Value = 1.063337300787233
$ws.Range("M13").Value = 1.075025774926894
$ws.Range("N13").Value = 1.024131990862558
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.056119050469896
$ws.Range("D14").Value = 1.063494167844377
$ws.Range("E14").Value = 1.060028907816437
$ws.Range("F14").Value = 1.071752977923611
$ws.Range("I14").Value = 1.055345332571332
$ws.Range("J14").Value = 1.062546041577768
$ws.Range("K14").Value = 1.066965862629562
$ws.Range("L14").Value = 1.06351278423135
$ws.Range("M14").Value = 1.075195989590552
$ws.Range("N14").Value = 1.024189325172384
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.056252149715929
$ws.Range("D15").Value = 1.063604263581384
$ws.Range("E15").Value = 1.060156448899491
$ws.Range("F15").Value = 1.071877014444702
$ws.Range("I15").Value = 1.055396450190838
$ws.Range("J15").Value = 1.062642280827399
$ws.Range("K15").Value = 1.067056651276694
$ws.Range("L15").Value = 1.063620883714026
$ws.Range("M15").Value = 1.075300843790255
$ws.Range("N15").Value = 1.024224627071955
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.057026597734205
$ws.Range("D16").Value = 1.064244865276529
$ws.Range("E16").Value = 1.060898675378744
$ws.Range("F16").Value = 1.072598862284603
$ws.Range("I16").Value = 1.055693392320731
$ws.Range("J16").Value = 1.06320200036783
$ws.Range("K16").Value = 1.067584651394353
$ws.Range("L16").Value = 1.064249740533208
$ws.Range("M16").Value = 1.075910826544655
$ws.Range("N16").Value = 1.024429738630749
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.057512173610056
$ws.Range("D17").Value = 1.06464652061587
$ws.Range("E17").Value = 1.061364152895166
$ws.Range("F17").Value = 1.073051573676586
$ws.Range("I17").Value = 1.055879142569748
$ws.Range("J17").Value = 1.063552716814978
$ws.Range("K17").Value = 1.067915475679052
$ws.Range("L17").Value = 1.064643918103464
$ws.Range("M17").Value = 1.076293177828473
$ws.Range("N17").Value = 1.024558082526794
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.05779532543777
$ws.Range("D18").Value = 1.064880736190356
$ws.Range("E18").Value = 1.061635622312937
$ws.Range("F18").Value = 1.073315602805775
$ws.Range("I18").Value = 1.055987302851832
$ws.Range("J18").Value = 1.063757147265603
$ws.Range("K18").Value = 1.06810830494532
$ws.Range("L18").Value = 1.064873732307644
$ws.Range("M18").Value = 1.076516098537895
$ws.Range("N18").Value = 1.024632829273053
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.05789186004198
$ws.Range("D19").Value = 1.064960586992262
$ws.Range("E19").Value = 1.06172818038596
$ws.Range("F19").Value = 1.073405624936916
$ws.Range("I19").Value = 1.05602415142598
$ws.Range("J19").Value = 1.063826829678067
$ws.Range("K19").Value = 1.068174031913967
$ws.Range("L19").Value = 1.064952075630683
$ws.Range("M19").Value = 1.076592092126676
$ws.Range("N19").Value = 1.024658296653392
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.057460083838465
$ws.Range("D20").Value = 1.064603433351074
$ws.Range("E20").Value = 1.061314215280195
$ws.Range("F20").Value = 1.073003005125557
$ws.Range("I20").Value = 1.055859232440272
$ws.Range("J20").Value = 1.063515102397489
$ws.Range("K20").Value = 1.067879995337363
$ws.Range("L20").Value = 1.064601637242971
$ws.Range("M20").Value = 1.076252165341298
$ws.Range("N20").Value = 1.024544324246327
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.056055432512997
$ws.Range("D21").Value = 1.063441544994707
$ws.Range("E21").Value = 1.059967948624829
$ws.Range("F21").Value = 1.071693694051575
$ws.Range("I21").Value = 1.055320891040767
$ws.Range("J21").Value = 1.06250003721972
$ws.Range("K21").Value = 1.06692246343544
$ws.Range("L21").Value = 1.063461113247294
$ws.Range("M21").Value = 1.075145869925069
$ws.Range("N21").Value = 1.02417244657756
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.055171561312374
$ws.Range("D22").Value = 1.062710436196542
$ws.Range("E22").Value = 1.059121157089725
$ws.Range("F22").Value = 1.070870193920828
$ws.Range("I22").Value = 1.054980741190152
$ws.Range("J22").Value = 1.06186057998899
$ws.Range("K22").Value = 1.066319195537669
$ws.Range("L22").Value = 1.062743076822764
$ws.Range("M22").Value = 1.074449397249547
$ws.Range("N22").Value = 1.023937599798061
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.055640186694881
$ws.Range("D23").Value = 1.063098066799873
$ws.Range("E23").Value = 1.059570090443019
$ws.Range("F23").Value = 1.071306774652145
$ws.Range("I23").Value = 1.055161220119059
$ws.Range("J23").Value = 1.06219968717748
$ws.Range("K23").Value = 1.066639116490052
$ws.Range("L23").Value = 1.063123811780164
$ws.Range("M23").Value = 1.074818696711376
$ws.Range("N23").Value = 1.024062194624731
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.057483621199954
$ws.Range("D24").Value = 1.064622902826639
$ws.Range("E24").Value = 1.061336780049295
$ws.Range("F24").Value = 1.073024951255155
$ws.Range("I24").Value = 1.055868229541429
$ws.Range("J24").Value = 1.063532099155462
$ws.Range("K24").Value = 1.067896027793483
$ws.Range("L24").Value = 1.064620742461811
$ws.Range("M24").Value = 1.076270697422706
$ws.Range("N24").Value = 1.024550541373555
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.059619840371451
$ws.Range("D25").Value = 1.066389918449394
$ws.Range("E25").Value = 1.063385532063668
$ws.Range("F25").Value = 1.075017639951743
$ws.Range("I25").Value = 1.056681434063752
$ws.Range("J25").Value = 1.065072941527023
$ws.Range("K25").Value = 1.069349316859035
$ws.Range("L25").Value = 1.066353818860635
$ws.Range("M25").Value = 1.077951813834221
$ws.Range("N25").Value = 1.025112767981979
